$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "Discrete random variables"
$ws.Range("D9").Value = "Continuous random variables"

$ws.Range("D10").Select()
